$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgals1"
$ws.Range("C2").Value = "Ptprc"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 37.833119
$ws.Range("H2").Value = 113.499357
$ws.Range("I2").Value = 0.2771305381131279
$ws.Range("J2").Value = 0.2771305381131279
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 243.7171333333333
$ws.Range("N2").Value = 731.1514
$ws.Range("O2").Value = 0.9993032963424349
$ws.Range("P2").Value = 0.999303296342435
$ws.Range("Q2").Value = 9220.579307738868
$ws.Range("R2").Value = 82985.2137696498
$ws.Range("S2").Value = 0.2769374602536015
$ws.Range("T2").Value = 0.2769374602536016

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgals1"
$ws.Range("C3").Value = "Ptprc"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 37.833119
$ws.Range("H3").Value = 113.499357
$ws.Range("I3").Value = 0.2771305381131279
$ws.Range("J3").Value = 0.2771305381131279
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1280236666666667
$ws.Range("N3").Value = 0.384071
$ws.Range("O3").Value = 0.0005249301530839377
$ws.Range("P3").Value = 0.0005249301530839377
$ws.Range("Q3").Value = 4.843534615816334
$ws.Range("R3").Value = 43.591811542347
$ws.Range("S3").Value = 0.0001454741757959583
$ws.Range("T3").Value = 0.0001454741757959583

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lgals1"
$ws.Range("C4").Value = "Ptprc"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 37.833119
$ws.Range("H4").Value = 113.499357
$ws.Range("I4").Value = 0.2771305381131279
$ws.Range("J4").Value = 0.2771305381131279
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04189333333333333
$ws.Range("N4").Value = 0.12568
$ws.Range("O4").Value = 0.0001717735044811748
$ws.Range("P4").Value = 0.0001717735044811748
$ws.Range("Q4").Value = 1.584955465306667
$ws.Range("R4").Value = 14.26459918776
$ws.Range("S4").Value = 0.00004760368373044576
$ws.Range("T4").Value = 0.00004760368373044576

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgals1"
$ws.Range("C5").Value = "Ptprc"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.88336466666667
$ws.Range("H5").Value = 191.650094
$ws.Range("I5").Value = 0.4679506129682439
$ws.Range("J5").Value = 0.467950612968244
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 243.7171333333333
$ws.Range("N5").Value = 731.1514
$ws.Range("O5").Value = 0.9993032963424349
$ws.Range("P5").Value = 0.999303296342435
$ws.Range("Q5").Value = 15569.47050424796
$ws.Range("R5").Value = 140125.2345382316
$ws.Range("S5").Value = 0.4676245900646291
$ws.Range("T5").Value = 0.4676245900646292

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lgals1"
$ws.Range("C6").Value = "Ptprc"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.88336466666667
$ws.Range("H6").Value = 191.650094
$ws.Range("I6").Value = 0.4679506129682439
$ws.Range("J6").Value = 0.467950612968244
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1280236666666667
$ws.Range("N6").Value = 0.384071
$ws.Range("O6").Value = 0.0005249301530839377
$ws.Range("P6").Value = 0.0005249301530839377
$ws.Range("Q6").Value = 8.178582583630446
$ws.Range("R6").Value = 73.60724325267401
$ws.Range("S6").Value = 0.0002456413869011428
$ws.Range("T6").Value = 0.0002456413869011428

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lgals1"
$ws.Range("C7").Value = "Ptprc"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.88336466666667
$ws.Range("H7").Value = 191.650094
$ws.Range("I7").Value = 0.4679506129682439
$ws.Range("J7").Value = 0.467950612968244
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.04189333333333333
$ws.Range("N7").Value = 0.12568
$ws.Range("O7").Value = 0.0001717735044811748
$ws.Range("P7").Value = 0.0001717735044811748
$ws.Range("Q7").Value = 2.676287090435556
$ws.Range("R7").Value = 24.08658381392
$ws.Range("S7").Value = 0.00008038151671366913
$ws.Range("T7").Value = 0.00008038151671366914

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lgals1"
$ws.Range("C8").Value = "Ptprc"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 34.80083866666666
$ws.Range("H8").Value = 104.402516
$ws.Range("I8").Value = 0.2549188489186281
$ws.Range("J8").Value = 0.2549188489186282
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 243.7171333333333
$ws.Range("N8").Value = 731.1514
$ws.Range("O8").Value = 0.9993032963424349
$ws.Range("P8").Value = 0.999303296342435
$ws.Range("Q8").Value = 8481.560637435821
$ws.Range("R8").Value = 76334.04573692239
$ws.Range("S8").Value = 0.2547412460242042
$ws.Range("T8").Value = 0.2547412460242043

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lgals1"
$ws.Range("C9").Value = "Ptprc"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 34.80083866666666
$ws.Range("H9").Value = 104.402516
$ws.Range("I9").Value = 0.2549188489186281
$ws.Range("J9").Value = 0.2549188489186282
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1280236666666667
$ws.Range("N9").Value = 0.384071
$ws.Range("O9").Value = 0.0005249301530839377
$ws.Range("P9").Value = 0.0005249301530839377
$ws.Range("Q9").Value = 4.455330969181778
$ws.Range("R9").Value = 40.09797872263599
$ws.Range("S9").Value = 0.0001338145903868367
$ws.Range("T9").Value = 0.0001338145903868367

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Lgals1"
$ws.Range("C10").Value = "Ptprc"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 34.80083866666666
$ws.Range("H10").Value = 104.402516
$ws.Range("I10").Value = 0.2549188489186281
$ws.Range("J10").Value = 0.2549188489186282
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.04189333333333333
$ws.Range("N10").Value = 0.12568
$ws.Range("O10").Value = 0.0001717735044811748
$ws.Range("P10").Value = 0.0001717735044811748
$ws.Range("Q10").Value = 1.457923134542222
$ws.Range("R10").Value = 13.12130821088
$ws.Range("S10").Value = 0.00004378830403705989
$ws.Range("T10").Value = 0.0000437883040370599
